# Generate Report for Handback
# This script updates the localization-status workbook to reflect that
# the zh-cn and de-de handback packages have been generated:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The "Latest Target File" / "Latest Handback File" columns on the per
#    -language sheets are populated with links/file names
#  - The "Latest Handback DateTime" on the de-de sheet is stamped with the
#    generation time

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Update the Status text everywhere it appears ("Ready for handoff"
#    is a shared string used on all three sheets).
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (I) / Latest Handback File (J)
# ---------------------------------------------------------------------
$null = $ws2.Hyperlinks.Add(
    $ws2.Cells.Item(2, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8608879a795b3c6a0467b24ac74131b130379cf3/e2e/8b18e12a-bc0f-4f4d-9fab-6873b5687b2b.md",
    $null,
    $null,
    "8b18e12a-bc0f-4f4d-9fab-6873b5687b2b.md")
$ws2.Range("J2").Value = "8b18e12a-bc0f-4f4d-9fab-6873b5687b2b.0c48187a2c4c5146b3ef183b452de327c6ab8cd8.zh-cn.xlf"

$null = $ws2.Hyperlinks.Add(
    $ws2.Cells.Item(3, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8608879a795b3c6a0467b24ac74131b130379cf3/e2e/d3e92616-89a4-46c4-bc1f-22960e5429c6.md",
    $null,
    $null,
    "d3e92616-89a4-46c4-bc1f-22960e5429c6.md")
$ws2.Range("J3").Value = "d3e92616-89a4-46c4-bc1f-22960e5429c6.54fb73fde87580050e96006e8e282e749052cfcf.zh-cn.xlf"

# ---------------------------------------------------------------------
# 3. de-de sheet: fill in Latest Target File (I) / Latest Handback File (J)
#    / Latest Handback DateTime (K)
# ---------------------------------------------------------------------
$null = $ws3.Hyperlinks.Add(
    $ws3.Cells.Item(2, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8608879a795b3c6a0467b24ac74131b130379cf3/e2e/8b18e12a-bc0f-4f4d-9fab-6873b5687b2b.md",
    $null,
    $null,
    "8b18e12a-bc0f-4f4d-9fab-6873b5687b2b.md")
$ws3.Range("J2").Value = "8b18e12a-bc0f-4f4d-9fab-6873b5687b2b.0c48187a2c4c5146b3ef183b452de327c6ab8cd8.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-02 04:31:36"

$null = $ws3.Hyperlinks.Add(
    $ws3.Cells.Item(3, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8608879a795b3c6a0467b24ac74131b130379cf3/e2e/d3e92616-89a4-46c4-bc1f-22960e5429c6.md",
    $null,
    $null,
    "d3e92616-89a4-46c4-bc1f-22960e5429c6.md")
$ws3.Range("J3").Value = "d3e92616-89a4-46c4-bc1f-22960e5429c6.54fb73fde87580050e96006e8e282e749052cfcf.de-de.xlf"
$ws3.Range("K3").Value = "2016-09-02 04:31:36"

# ---------------------------------------------------------------------
# 4. zh-cn Latest Handback DateTime (K) gets stamped via the shared
#    "0001-01-01 00:00:00" -> "2016-09-02 04:31:28" text replacement.
# ---------------------------------------------------------------------
$ws2.Range("K2").Value = "2016-09-02 04:31:28"
$ws2.Range("K3").Value = "2016-09-02 04:31:28"

# ---------------------------------------------------------------------
# 5. Widen columns to fit the newly-written long values (mirrors the
#    report generator's auto-fit/column sizing pass).
# ---------------------------------------------------------------------
$wideWidth = 29.17      # renders to ~29.98-30 character units
$fullWidth = 39.166666666666664   # renders to exactly 40 character units

$ws1.Columns.Item(5).ColumnWidth = $wideWidth   # zh-cn status column
$ws1.Columns.Item(6).ColumnWidth = $wideWidth   # de-de status column

$ws2.Columns.Item(3).ColumnWidth = $wideWidth   # Status
$ws2.Columns.Item(9).ColumnWidth = $fullWidth   # Latest Target File
$ws2.Columns.Item(10).ColumnWidth = $fullWidth  # Latest Handback File

$ws3.Columns.Item(3).ColumnWidth = $wideWidth   # Status
$ws3.Columns.Item(9).ColumnWidth = $fullWidth   # Latest Target File
$ws3.Columns.Item(10).ColumnWidth = $fullWidth  # Latest Handback File
